# Drop in RMI script files
# - Remove the "Texas Notes" worksheet entirely.
# - Replace the formula in PDiCECpDoC!B2 (which referenced 'Texas Notes'!B10)
#   with the plain literal value it used to resolve to (0.13).

$wb = $excel.ActiveWorkbook

# Update the PDiCECpDoC sheet's B2 cell to a hard-coded value before removing
# the sheet it depended on.
$target = $wb.Worksheets.Item("PDiCECpDoC")
$target.Range("B2").Value = 0.13

# Remove the now-unused "Texas Notes" sheet.
$notes = $wb.Worksheets.Item("Texas Notes")
$notes.Delete()

# Make sure the "About" sheet is the active / selected sheet, matching the
# final saved state.
$about = $wb.Worksheets.Item("About")
$about.Activate()
